$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as Text so numeric-looking strings
# (e.g. "1.007") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.045.55"
$ws.Range("E2").Value = "  -1.30%  "
$ws.Range("D3").Value = "1.841.50"
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").Value = "313.15"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").Value = "0.4630"
$ws.Range("E7").Value = "  -0.71%  "
$ws.Range("D8").Value = "0.3686"
$ws.Range("E8").Value = "  -1.10%  "
$ws.Range("D9").Value = "0.07252"
$ws.Range("E9").Value = "  -1.98%  "
$ws.Range("D10").Value = "0.8830"
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "0.07824"
$ws.Range("E11").Value = "  -1.72%  "
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "19.76"
$ws.Range("E12").Value = "  -2.04%  "
$ws.Range("D13").Value = "1.867.13"
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").Value = "5.381"
$ws.Range("E14").Value = "  -1.06%  "
$ws.Range("D15").Value = "6.486"
$ws.Range("E15").Value = "  -2.21%  "
$ws.Range("D16").Value = "91.14"
$ws.Range("E16").Value = "  -1.68%  "
$ws.Range("D17").Value = "1.007"
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("D18").Value = "0.000008818"
$ws.Range("E18").Value = "  -1.55%  "
$ws.Range("D19").Value = "1.004"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").Value = "27.111.80"
$ws.Range("E20").Value = "  -1.16%  "
$ws.Range("D21").Value = "14.60"
$ws.Range("E21").Value = "  -2.26%  "
$ws.Range("D22").Value = "5.030"
$ws.Range("E22").Value = "  -2.63%  "
$ws.Range("D23").Value = "10.52"
$ws.Range("E23").Value = "  -0.83%  "
$ws.Range("D24").Value = "2.028"
$ws.Range("E24").Value = "  +8.82%  "
$ws.Range("D25").Value = "150.84"
$ws.Range("E25").Value = "  -1.12%  "
$ws.Range("D26").Value = "18.32"
$ws.Range("E26").Value = "  -1.49%  "
$ws.Range("D27").Value = "2.020"
$ws.Range("E27").Value = "  -3.58%  "
$ws.Range("D28").Value = "115.60"
$ws.Range("E28").Value = "  -1.55%  "
$ws.Range("D29").Value = "4.985"
$ws.Range("E29").Value = "  -3.55%  "
$ws.Range("D30").Value = "0.08858"
$ws.Range("E30").Value = "  -0.62%  "
$ws.Range("D31").Value = "3.148"
$ws.Range("E31").Value = "  +6.16%  "
$ws.Range("D32").Value = "0.7746"
$ws.Range("E32").Value = "  +2.53%  "
$ws.Range("D33").Value = "4.494"
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("D34").Value = "1.146"
$ws.Range("E34").Value = "  -1.17%  "
$ws.Range("D35").Value = "2.675"
$ws.Range("E35").Value = "  +3.01%  "
$ws.Range("D36").Value = "1.100"
$ws.Range("E36").Value = "  +1.32%  "
$ws.Range("D37").Value = "0.01937"
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("D38").Value = "0.05200"
$ws.Range("E38").Value = "  -1.92%  "
$ws.Range("D39").Value = "2.944"
$ws.Range("E39").Value = "  -1.56%  "
$ws.Range("D40").Value = "7.003"
$ws.Range("E40").Value = "  -2.36%  "
$ws.Range("D41").Value = "0.5028"
$ws.Range("E41").Value = "  -3.64%  "
$ws.Range("D42").Value = "0.1609"
$ws.Range("E42").Value = "  -2.27%  "
$ws.Range("D43").Value = "8.447"
$ws.Range("E43").Value = "  +0.98%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "10.39"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.4722"
$ws.Range("E45").Value = "  -3.92%  "
$ws.Range("D46").Value = "1.003"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").Value = "103.05"
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("D48").Value = "1.626"
$ws.Range("E48").Value = "  -1.21%  "
$ws.Range("D49").Value = "0.06171"
$ws.Range("E49").Value = "  -1.67%  "
$ws.Range("D50").Value = "65.39"
$ws.Range("E50").Value = "  -0.40%  "
$ws.Range("D51").Value = "36.36"
$ws.Range("E51").Value = "  -2.33%  "

# Restore default (no explicit number format) now that values are set.
$ws.Range("D2:D51").ClearFormats()

